# Fine tuning ecYali: f = 0.5 and kcat curation
#
# 1. Rename the only sheet ("Planilha1") to "F 0.3".
# 2. Correct a couple of kcat/notes values on "F 0.3".
# 3. Append the new DLKcat-priority-3 curated rows (41-45) to "F 0.3".
# 4. Add a brand-new sheet "F 0.5" (inserted after "F 0.3") holding the
#    f=0.5 sensitivity/DLKcat-priority-2 tuned subset, and make it the
#    active tab/selected sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "F 0.3"

# --- corrections on existing rows -----------------------------------------

# Row 14: kcat bumped to 891 and note replaced with the sensitivityTuning
# explanation (multiplying highest kcat/Km * Km of E. coli EC 2.7.1.29).
$ws1.Cells.Item(14, 4).Value = 891
$ws1.Cells.Item(14, 6).Value = "Output of sensitivityTuning. Calculated by multiplying highest kcat/Km * Km of Escherichia coli (EC 2.7.1.29)."

# Row 18: typo fix only ("Brenda EC" -> "EC") in the note text.
$ws1.Cells.Item(18, 6).Value = "Limits model with DLKcat priority 5. Calculatedby multiplying highest kcat/Km * Km (EC 1.14.19.41)"

# --- newly curated rows appended to "F 0.3" --------------------------------

$newRows = @(
    @("Q6CGV2", "YALI0A15950g", "YALI0A15950g", 507.794,  "y000910", "Limits model with DLKcat priority 3. Calculated from the specific activity of S. cerevisiae (EC 3.6.1.31)", 1),
    @("Q6C7Y2", "YALI0D24409g", "YALI0D24409g", 4.97,     "y000912", "Limits model with DLKcat priority 3. Using kcat of E. coli (EC 2.1.2.3)", 1),
    @("Q99148", "YALI0F21010g", "YALI0F21010g", 7.2,      "y000914", "Limits model with DLKcat priority 3. Using kcat of Gallus gallus (EC 6.3.4.13)", 1),
    @("P38997", "YALI0B15444g", "YALI0B15444g", 71.7479,  "y000988", "Limits model with DLKcat priority 3. Calculated from the specific activity of S. cerevisiae (EC 1.5.1.7)", 1),
    @("Q6CDK7", "YALI0B23188g", "YALI0B23188g", 40,       "y200001", "Limits model with DLKcat priority 3. Using kcat of Gallus gallus (EC 6.3.4.13)", 1)
)

$row = 41
foreach ($r in $newRows) {
    $ws1.Cells.Item($row, 1).Value = $r[0]
    $ws1.Cells.Item($row, 2).Value = $r[1]
    $ws1.Cells.Item($row, 3).Value = $r[2]
    $ws1.Cells.Item($row, 4).Value = $r[3]
    $ws1.Cells.Item($row, 5).Value = $r[4]
    $ws1.Cells.Item($row, 6).Value = $r[5]
    $ws1.Cells.Item($row, 7).Value = $r[6]
    $row = $row + 1
}

$ws1.Range("F14").Select()

# --- new "F 0.5" sheet ------------------------------------------------------

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "F 0.5"

$header = @("proteins", "genes", "gene_name", "kcat", "rxns", "notes", "stoicho")
for ($c = 1; $c -le $header.Length; $c++) {
    $ws2.Cells.Item(1, $c).Value = $header[$c - 1]
}

$sheet2Rows = @(
    @("Q6C791", "YALI0E02728g", "YALI0E02728g", 9.592,    "y000027", "Output of sensitivityTuning. Calculated from the specific activity of S. oneidensis (Brenda EC 4.2.1.117)", 1),
    @("Q6C120", "YALI0F19910g", "YALI0F19910g", 17.0475,  "y000029", "Limits model with DLKcat priority 2.  Calculated from the specific activity of N. crassa (Brenda)", 1),
    @("Q6CAF8", "YALI0D03135g", "YALI0D03135g", 3.81,     "y000096", "Output of sensitivityTuning. Got the highest value in Brenda for matching substrate (Meiothermus ruber).", 1),
    @("Q6C231", "YALI0F11297g", "YALI0F11297g", 10.1681,  "y000238, y000239, y000240", "Output of sensitivityTuning. Calculated from the specific activity of R. norvegicus (Brenda EC 1.14.18.9)", 1),
    @("Q6CAY2", "YALI0C23408g", "YALI0C23408g", 10.3331,  $null,     "Limits model with DLKcat priority 2.  Calculated from the specific activity of N. crassa (EC 4.2.1.9)", 1),
    @("Q6C564", "YALI0E20691g", "YALI0E20691g", 42.93,    "y000354", "Limits model with DLKcat priority 2. Calculated by multiplying highest kcat/Km * Km of Escherichia coli (EC 2.7.1.29)", 1),
    @("Q6C1F3", "YALI0F16819g", "YALI0F16819g", 230,      "y000366", "Limits model with DLKcat priority 2. Using kcat of S. cerevisiae (EC 4.2.1.11)", 1),
    @("Q6C6H1", "YALI0E09603g", "YALI0E09603g", 130.5365, "y000470", "Limits model with DLKcat priority 2.  Calculated from the specific activity of Aspergillus nidulans (EC 1.4.1.2)", 1),
    @("Q6CGV2", "YALI0A15950g", "YALI0A15950g", 507.794,  "y000910", "Limits model with DLKcat priority 2. Calculated from the specific activity of S. cerevisiae (EC 3.6.1.31)", 1),
    @("Q6C7Y2", "YALI0D24409g", "YALI0D24409g", 4.97,     "y000912", "Limits model with DLKcat priority 2. Using kcat of E. coli (EC 2.1.2.3)", 1),
    @("P38997", "YALI0B15444g", "YALI0B15444g", 71.7479,  "y000988", "Limits model with DLKcat priority 2. Calculated from the specific activity of S. cerevisiae (EC 1.5.1.7)", 1)
)

$row = 2
foreach ($r in $sheet2Rows) {
    $ws2.Cells.Item($row, 1).Value = $r[0]
    $ws2.Cells.Item($row, 2).Value = $r[1]
    $ws2.Cells.Item($row, 3).Value = $r[2]
    $ws2.Cells.Item($row, 4).Value = $r[3]
    if ($r[4] -ne $null) {
        $ws2.Cells.Item($row, 5).Value = $r[4]
    }
    $ws2.Cells.Item($row, 6).Value = $r[5]
    $ws2.Cells.Item($row, 7).Value = $r[6]
    $row = $row + 1
}

# Column widths approximating the autofit state captured on the real sheet.
$ws2.Columns.Item(2).ColumnWidth = 13.140625
$ws2.Columns.Item(3).ColumnWidth = 13.5703125
$ws2.Columns.Item(5).ColumnWidth = 24.85546875
$ws2.Columns.Item(6).ColumnWidth = 98.85546875

$ws2.Activate()
$ws2.Range("F14").Select()
